$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update time_taken (column F) timestamps on the "data" sheet ---
$ws.Cells.Item(2,6).Value = "2021-10-05 14:21:42.020198"
$ws.Cells.Item(3,6).Value = "2021-10-05 14:21:42.020207"
$ws.Cells.Item(4,6).Value = "2021-10-05 14:21:42.020211"
$ws.Cells.Item(5,6).Value = "2021-10-05 14:21:42.020214"
$ws.Cells.Item(6,6).Value = "2021-10-05 14:21:42.020217"
$ws.Cells.Item(7,6).Value = "2021-10-05 14:21:42.020219"
$ws.Cells.Item(8,6).Value = "2021-10-05 14:21:42.020222"
$ws.Cells.Item(9,6).Value = "2021-10-05 14:21:42.020224"
$ws.Cells.Item(10,6).Value = "2021-10-05 14:21:42.020227"
$ws.Cells.Item(11,6).Value = "2021-10-05 14:21:42.020230"
$ws.Cells.Item(12,6).Value = "2021-10-05 14:21:42.020232"
$ws.Cells.Item(13,6).Value = "2021-10-05 14:21:42.020235"
$ws.Cells.Item(14,6).Value = "2021-10-05 14:21:42.020238"
$ws.Cells.Item(15,6).Value = "2021-10-05 14:21:42.020240"
$ws.Cells.Item(16,6).Value = "2021-10-05 14:21:42.020243"
$ws.Cells.Item(17,6).Value = "2021-10-05 14:21:42.020245"
$ws.Cells.Item(18,6).Value = "2021-10-05 14:21:42.020248"
$ws.Cells.Item(19,6).Value = "2021-10-05 14:21:42.020251"
$ws.Cells.Item(20,6).Value = "2021-10-05 14:21:42.020253"
$ws.Cells.Item(21,6).Value = "2021-10-05 14:21:42.020256"
$ws.Cells.Item(22,6).Value = "2021-10-05 14:21:42.020258"
$ws.Cells.Item(23,6).Value = "2021-10-05 14:21:42.020261"
$ws.Cells.Item(24,6).Value = "2021-10-05 14:21:42.020264"
$ws.Cells.Item(25,6).Value = "2021-10-05 14:21:42.020266"
$ws.Cells.Item(26,6).Value = "2021-10-05 14:21:42.020269"
$ws.Cells.Item(27,6).Value = "2021-10-05 14:21:42.020272"
$ws.Cells.Item(28,6).Value = "2021-10-05 14:21:42.020274"
$ws.Cells.Item(29,6).Value = "2021-10-05 14:21:42.020277"
$ws.Cells.Item(30,6).Value = "2021-10-05 14:21:42.020279"
$ws.Cells.Item(31,6).Value = "2021-10-05 14:21:42.020282"
$ws.Cells.Item(32,6).Value = "2021-10-05 14:21:42.020285"
$ws.Cells.Item(33,6).Value = "2021-10-05 14:21:42.020289"
$ws.Cells.Item(34,6).Value = "2021-10-05 14:21:42.020292"
$ws.Cells.Item(35,6).Value = "2021-10-05 14:21:42.020295"
$ws.Cells.Item(36,6).Value = "2021-10-05 14:21:42.020298"
$ws.Cells.Item(37,6).Value = "2021-10-05 14:21:42.020300"
$ws.Cells.Item(38,6).Value = "2021-10-05 14:21:42.020303"
$ws.Cells.Item(39,6).Value = "2021-10-05 14:21:42.020305"
$ws.Cells.Item(40,6).Value = "2021-10-05 14:21:42.020308"
$ws.Cells.Item(41,6).Value = "2021-10-05 14:21:42.020310"
$ws.Cells.Item(42,6).Value = "2021-10-05 14:21:42.020313"
$ws.Cells.Item(43,6).Value = "2021-10-05 14:21:42.020316"
$ws.Cells.Item(44,6).Value = "2021-10-05 14:21:42.020318"
$ws.Cells.Item(45,6).Value = "2021-10-05 14:21:42.020321"
$ws.Cells.Item(46,6).Value = "2021-10-05 14:21:42.020323"
$ws.Cells.Item(47,6).Value = "2021-10-05 14:21:42.020326"
$ws.Cells.Item(48,6).Value = "2021-10-05 14:21:42.020328"
$ws.Cells.Item(49,6).Value = "2021-10-05 14:21:42.020331"
$ws.Cells.Item(50,6).Value = "2021-10-05 14:21:42.020333"
$ws.Cells.Item(51,6).Value = "2021-10-05 14:21:42.020336"
$ws.Cells.Item(52,6).Value = "2021-10-05 14:21:42.020339"
$ws.Cells.Item(53,6).Value = "2021-10-05 14:21:42.020341"
$ws.Cells.Item(54,6).Value = "2021-10-05 14:21:42.020344"
$ws.Cells.Item(55,6).Value = "2021-10-05 14:21:42.020347"
$ws.Cells.Item(56,6).Value = "2021-10-05 14:21:42.020350"
$ws.Cells.Item(57,6).Value = "2021-10-05 14:21:42.020352"
$ws.Cells.Item(58,6).Value = "2021-10-05 14:21:42.020355"
$ws.Cells.Item(59,6).Value = "2021-10-05 14:21:42.020358"
$ws.Cells.Item(60,6).Value = "2021-10-05 14:21:42.020360"
$ws.Cells.Item(61,6).Value = "2021-10-05 14:21:42.020363"
$ws.Cells.Item(62,6).Value = "2021-10-05 14:21:42.020365"
$ws.Cells.Item(63,6).Value = "2021-10-05 14:21:42.020368"
$ws.Cells.Item(64,6).Value = "2021-10-05 14:21:42.020371"
$ws.Cells.Item(65,6).Value = "2021-10-05 14:21:42.020373"
$ws.Cells.Item(66,6).Value = "2021-10-05 14:21:42.020377"
$ws.Cells.Item(67,6).Value = "2021-10-05 14:21:42.020380"
$ws.Cells.Item(68,6).Value = "2021-10-05 14:21:42.020383"
$ws.Cells.Item(69,6).Value = "2021-10-05 14:21:42.020385"
$ws.Cells.Item(70,6).Value = "2021-10-05 14:21:42.020388"
$ws.Cells.Item(71,6).Value = "2021-10-05 14:21:42.020391"
$ws.Cells.Item(72,6).Value = "2021-10-05 14:21:42.020394"
$ws.Cells.Item(73,6).Value = "2021-10-05 14:21:42.020396"
$ws.Cells.Item(74,6).Value = "2021-10-05 14:21:42.020399"
$ws.Cells.Item(75,6).Value = "2021-10-05 14:21:42.020402"
$ws.Cells.Item(76,6).Value = "2021-10-05 14:21:42.020404"
$ws.Cells.Item(77,6).Value = "2021-10-05 14:21:42.020407"
$ws.Cells.Item(78,6).Value = "2021-10-05 14:21:42.020412"
$ws.Cells.Item(79,6).Value = "2021-10-05 14:21:42.020415"
$ws.Cells.Item(80,6).Value = "2021-10-05 14:21:42.020418"
$ws.Cells.Item(81,6).Value = "2021-10-05 14:21:42.020421"
$ws.Cells.Item(82,6).Value = "2021-10-05 14:21:42.020423"
$ws.Cells.Item(83,6).Value = "2021-10-05 14:21:42.020426"
$ws.Cells.Item(84,6).Value = "2021-10-05 14:21:42.020428"
$ws.Cells.Item(85,6).Value = "2021-10-05 14:21:42.020431"
$ws.Cells.Item(86,6).Value = "2021-10-05 14:21:42.020434"
$ws.Cells.Item(87,6).Value = "2021-10-05 14:21:42.020436"
$ws.Cells.Item(88,6).Value = "2021-10-05 14:21:42.020439"
$ws.Cells.Item(89,6).Value = "2021-10-05 14:21:42.020441"
$ws.Cells.Item(90,6).Value = "2021-10-05 14:21:42.020444"
$ws.Cells.Item(91,6).Value = "2021-10-05 14:21:42.020447"
$ws.Cells.Item(92,6).Value = "2021-10-05 14:21:42.020449"
$ws.Cells.Item(93,6).Value = "2021-10-05 14:21:42.020452"
$ws.Cells.Item(94,6).Value = "2021-10-05 14:21:42.020456"
$ws.Cells.Item(95,6).Value = "2021-10-05 14:21:42.020459"

# --- Add the "metadata" sheet after "data" ---
$newws = $wb.Worksheets.Add($null, $ws)
$newws.Name = "metadata"

# Header row (row 1), columns B:G
$newws.Cells.Item(1,2).Value = "data_name"
$newws.Cells.Item(1,3).Value = "data_id"
$newws.Cells.Item(1,4).Value = "data_version"
$newws.Cells.Item(1,5).Value = "data_version_created"
$newws.Cells.Item(1,6).Value = "panel_query_time"
$newws.Cells.Item(1,7).Value = "panel_get_request"

# Copy the header style from the "data" sheet's header row so we reuse
# the existing bold/bordered/centered style instead of creating a new one.
$ws.Cells.Item(1,2).Copy()
$newws.Range("B1:G1").PasteSpecial(-4122)

# Data row (row 2)
$newws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,1).Copy()
$newws.Cells.Item(2,1).PasteSpecial(-4122)

$newws.Cells.Item(2,2).Value = "Neonatal cholestasis"
$newws.Cells.Item(2,3).Value = 385

# Force "1.20" to stay text (not get parsed into the number 1.2)
$newws.Cells.Item(2,4).NumberFormat = "@"
$newws.Cells.Item(2,4).Value = "1.20"
$newws.Cells.Item(2,4).Style = "Normal"

$newws.Cells.Item(2,5).Value = "2021-09-06T10:15:19.555779Z"
$newws.Cells.Item(2,6).Value = "2021-10-05 14:21:42.016902"
$newws.Cells.Item(2,7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/385/?format=json"
